$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.297.30'
$ws.Range("E2").Value = '  +4.13%  '
$ws.Range("D3").Value = '1.713.41'
$ws.Range("E3").Value = '  +1.68%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +4.89%  '
$ws.Range("E9").Value = '  +2.31%  '
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '1.957.72'
$ws.Range("E12").Value = '  +1.78%  '
$ws.Range("D13").Value = '1.717.03'
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("E17").Value = '  +7.72%  '
$ws.Range("D18").Value = '28.271.56'
$ws.Range("E18").Value = '  +4.01%  '
$ws.Range("D19").Value = '0.0₃0747'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("E24").Value = '  -1.52%  '
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("E31").Value = '  +2.55%  '
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("D33").Value = '1.482.81'
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("E36").Value = '  +2.12%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -1.20%  '
$ws.Range("E39").Value = '  +0.71%  '
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").Value = '1.862.55'
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("E47").Value = '  +8.12%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0112'
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("E51").Value = '  -3.51%  '

# Force numeric-looking price values to remain text (inline string), matching source data
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.28'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.29'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.266'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0632'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.21'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.561'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.39'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '254.50'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.78'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.57'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.58'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.13'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.39'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.62'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0511'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.20'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.967'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.59'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.801'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.75'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '90.01'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.02'

# Reset style back to default (no explicit format) to match original workbook styling
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
